# EvidenciasMuestreoAprobados.xlsx - update saved view/selection and swap the
# last two header labels on row 1 (N1 <-> O1), matching the authored commit:
# "Se crean apis y funciones para la descarga de excel en modulo de
#  aprobados, rechazados y eventualidades..."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the header text in N1 and O1 ---------------------------------
# N1 was "FECHA DE VALIDACIÓN DE CONAGUA", O1 was "PORCENTAJE PAGO".
# After the edit they trade places: N1 -> "PORCENTAJE PAGO",
# O1 -> "FECHA DE VALIDACIÓN DE CONAGUA".
$n1Value = $ws.Range("N1").Value2
$o1Value = $ws.Range("O1").Value2

$ws.Range("N1").Value2 = $o1Value
$ws.Range("O1").Value2 = $n1Value

# --- Move the selection from L14 to M5 ----------------------------------
$ws.Range("M5").Select()
